$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for "even_MAG-GUT91631.fa" (row 5) was removed from the sheet;
# subsequent rows (GUT91672, GUT91675) shift up to fill the gap.
$ws.Rows("5:5").Delete()
